# Atualização de bases das ligas, do dia: 29-02-2024 às 07:50
#
# For each of the row pairs below, the two rows describe the same matchday
# but the Home/Away fixture (and all stats columns B, F..AC) got swapped
# between the two rows (column A "id" and the shared C/D/E columns stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns to swap between each row pair: B, and F through AC (everything
# except A, C, D, E which remain the same for the pair).
$colsToSwap = @(2) + (6..29)

$rowPairs = @(
    @(8, 9),
    @(10, 11),
    @(13, 14),
    @(18, 19),
    @(26, 27),
    @(32, 33),
    @(36, 37),
    @(47, 48),
    @(51, 52)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $colsToSwap) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value()
        $v2 = $cell2.Value()

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
